# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" worksheet (directly after "总计", before "2022-Q2")
# containing the Q3-2022 fund snapshot, and records the new quarter in the
# "总计" (summary) worksheet as its first data row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q3" worksheet by duplicating "2022-Q2" (so it inherits
#    identical layout/column widths/styles), inserted right before it.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($template)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Use a never-styled cell (H2, plain default style) as a formatting donor so
# that after we write the new figures as text the cells end up back on the
# workbook's default (unstyled) cell format - exactly like the other quarter
# sheets' D/E/F/G data cells.
$q3.Range("H2").Copy()
$q3.Range("D2:G3").PasteSpecial(-4122)

$valueCells = @("D2", "E2", "F2", "G2", "D3", "E3", "F3", "G3")
foreach ($addr in $valueCells) {
    $q3.Range($addr).NumberFormat = "@"
}

# Row 2: fund 008763 "天弘越南市场股票（QDII）A"
$q3.Range("D2").Value = "20.44"
$q3.Range("E2").Value = "90.19"
$q3.Range("F2").Value = "6.38"
$q3.Range("G2").Value = "1.3041"

# Row 3: fund 008764 "天弘越南市场股票（QDII）C"
$q3.Range("D3").Value = "15.02"
$q3.Range("E3").Value = "90.19"
$q3.Range("F3").Value = "6.38"
$q3.Range("G3").Value = "0.9583"

# Drop back to the unstyled default format (mirrors the source data's lack
# of an explicit style on these cells).
$q3.Range("H2").Copy()
$q3.Range("D2:G3").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Record the new quarter as the first data row of the "总计" sheet,
#    pushing the existing quarters down by one row and renumbering the
#    running index in column A.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Rows.Item(2).Insert(-4121)

# Re-apply the original data-row formatting (lost by Insert) to the new row,
# copying it from the row directly below (still formatted like a data row).
$summary.Range("A3:D3").Copy()
$summary.Range("A2:D2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 2
$summary.Range("D2").Value = 2.26

# Renumber the running index (column A) for the rows that shifted down.
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4

# ---------------------------------------------------------------------------
# 3) Restore the originally-selected tab ("2021-Q3", the last sheet).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()
